$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CCDeferredPlanCorp")
$ws.Range("B2").Value = "Mon Nov 17 03:23:16 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:25:31 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPC")
$ws.Range("B2").Value = "Mon Nov 17 03:28:21 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:29:22 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPS")
$ws.Range("B2").Value = "Mon Nov 17 03:30:17 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:31:14 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanCredit")
$ws.Range("B2").Value = "Mon Nov 17 03:26:28 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:27:24 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCorp")
$ws.Range("B2").Value = "Mon Nov 17 03:12:43 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:13:46 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 03:14:50 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:15:57 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Mon Nov 17 03:16:58 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:19:03 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPS")
$ws.Range("B2").Value = "Mon Nov 17 03:20:01 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:22:13 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageDataCC")
$ws.Range("B2").Value = "Mon Nov 17 04:59:08 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 05:01:28 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCorp")
$ws.Range("B2").Value = "Mon Nov 17 02:50:36 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:51:31 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCorp")
$ws.Range("B2").Value = "Mon Nov 17 02:34:43 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:35:56 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCorp")
$ws.Range("B2").Value = "Mon Nov 17 02:58:22 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:59:20 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPS")
$ws.Range("B2").Value = "Mon Nov 17 02:54:13 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:55:19 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPS")
$ws.Range("B2").Value = "Mon Nov 17 03:02:20 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:03:25 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCredit")
$ws.Range("B2").Value = "Mon Nov 17 02:32:44 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:33:43 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Mon Nov 17 02:36:59 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:37:41 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPS")
$ws.Range("B2").Value = "Mon Nov 17 02:38:44 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:39:52 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCredit")
$ws.Range("B2").Value = "Mon Nov 17 02:48:44 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:49:41 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCredit")
$ws.Range("B2").Value = "Mon Nov 17 02:56:18 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:57:22 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPC")
$ws.Range("B2").Value = "Mon Nov 17 02:52:21 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 02:53:22 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPC")
$ws.Range("B2").Value = "Mon Nov 17 03:00:21 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:01:19 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 03:48:15 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Mon Nov 17 03:49:17 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 03:44:18 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Mon Nov 17 03:45:22 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 03:46:14 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:47:14 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 03:41:56 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Mon Nov 17 03:43:09 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCorp")
$ws.Range("B2").Value = "Sun Nov 16 23:34:41 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:34:27 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCredit")
$ws.Range("B2").Value = "Mon Nov 17 03:35:35 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:36:23 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPC")
$ws.Range("B2").Value = "Mon Nov 17 03:37:24 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:38:20 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 03:39:37 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 03:40:46 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayCorp")
$ws.Range("B2").Value = "Mon Nov 17 04:42:34 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:43:35 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayCC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Mon Nov 17 04:39:30 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:41:47 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayPC")
$ws.Range("B2").Value = "Mon Nov 17 04:44:29 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:45:16 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayCorp")
$ws.Range("B2").Value = "Mon Nov 17 19:26:45 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 18:55:55 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 04:47:43 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:48:50 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayPS")
$ws.Range("B2").Value = "Mon Nov 17 04:56:41 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:57:50 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayPC")
$ws.Range("B2").Value = "Mon Nov 17 18:48:40 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 18:50:44 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Mon Nov 17 04:31:38 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:32:45 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayCC")
$ws.Range("B2").Value = "Mon Nov 17 19:36:05 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 19:38:12 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayPC")
$ws.Range("B2").Value = "Mon Nov 17 04:34:02 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:35:06 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Mon Nov 17 04:36:11 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:38:23 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Mon Nov 17 04:46:06 IST 2025"
$ws.Range("B3").Value = "Mon Nov 17 04:46:52 IST 2025"
